$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.902.65'
$ws.Range('E2').Value = '  -3.23%  '
$ws.Range('D3').Value = '2.275.42'
$ws.Range('E3').Value = '  -4.03%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '532.17'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -4.39%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '131.11'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.54%  '
$ws.Range('E7').Value = '  +0.35%  '
$ws.Range('E8').Value = '  +0.85%  '
$ws.Range('D9').Value = '2.271.99'
$ws.Range('E9').Value = '  -4.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0991'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -5.59%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.42'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.23%  '
$ws.Range('E12').Value = '  -0.30%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.329'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.56%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.45'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.65%  '
$ws.Range('D15').Value = '2.683.15'
$ws.Range('E15').Value = '  -4.05%  '
$ws.Range('D16').Value = '57.866.95'
$ws.Range('E16').Value = '  -3.25%  '
$ws.Range('E17').Value = '  -4.32%  '
$ws.Range('D18').Value = '2.262.50'
$ws.Range('E18').Value = '  -4.79%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.51'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -5.06%  '
$ws.Range('E20').Value = '  -5.71%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '311.92'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.55%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.35'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.42%  '
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '62.46'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.81%  '
$ws.Range('E25').Value = '  -3.42%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.14%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.97'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.31%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.27'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -7.70%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '170.34'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.02%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.71'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.05%  '
$ws.Range('D31').Value = '0.0₃0716'
$ws.Range('E31').Value = '  -5.44%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.74'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.31%  '
$ws.Range('E33').Value = '  -6.36%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.378'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.34%  '
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.70'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.21%  '
$ws.Range('E37').Value = '  -0.20%  '
$ws.Range('E38').Value = '  -7.20%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.88'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -6.01%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.48'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -6.39%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '141.45'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.46%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '285.64'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -10.61%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.42'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.01%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0950'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.63%  '
$ws.Range('E46').Value = '  -3.09%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.549'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.17%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '18.03'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -7.74%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0209'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.50%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '10.93'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.21%  '
$ws.Range('E51').Value = '  -0.60%  '
